# Generate Report for Handoff
# Update the status/handoff-date for the "ee06cae0-c551-4ad6-99fb-a3c8ada45c1d.md"
# file to reflect that it is now ready for a new handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the ee06cae0-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-31-11 16:31:58"

# --- zh-cn sheet: row 3 is the ee06cae0-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-11 16:31:55"

# --- de-de sheet: row 3 is the ee06cae0-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-11 16:31:58"
